# previsao_retorno.xlsx refresh: "atualizacao dos dados bibi e add"
#
# The source pipeline regenerated this report a bit later than the previous
# run, so every "meses sem comprar" (months-without-purchase) counter that
# was still INATIVO crept forward by ~0.1 month, and two clients (rows 71 and
# 111) picked up brand new purchase activity, updating their probability /
# purchase-count / date / situation fields (row 71 even flipped from INATIVO
# back to ATIVO).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple "situacao" (INATIVO - X.Y meses sem comprar) text bumps ---
$ws.Range("J5").Value   = "INATIVO - 14.0 meses sem comprar"
$ws.Range("J6").Value   = "INATIVO - 15.5 meses sem comprar"
$ws.Range("J15").Value  = "INATIVO - 39.0 meses sem comprar"
$ws.Range("J23").Value  = "INATIVO - 37.0 meses sem comprar"
$ws.Range("J37").Value  = "INATIVO - 31.4 meses sem comprar"
$ws.Range("J44").Value  = "INATIVO - 15.1 meses sem comprar"
$ws.Range("J48").Value  = "INATIVO - 6.7 meses sem comprar"
$ws.Range("J63").Value  = "INATIVO - 27.0 meses sem comprar"
$ws.Range("J74").Value  = "INATIVO - 6.9 meses sem comprar"
$ws.Range("J85").Value  = "INATIVO - 14.2 meses sem comprar"
$ws.Range("J93").Value  = "INATIVO - 15.6 meses sem comprar"
$ws.Range("J95").Value  = "INATIVO - 32.1 meses sem comprar"
$ws.Range("J99").Value  = "INATIVO - 36.4 meses sem comprar"
$ws.Range("J100").Value = "INATIVO - 6.7 meses sem comprar"
$ws.Range("J101").Value = "INATIVO - 13.7 meses sem comprar"

# --- Row 71 (id_cliente 19765, INGRID MORAES FERNANDES): new purchase ---
# moved the client from INATIVO back into ATIVO.
$ws.Range("B71").Value = 0.42
$ws.Range("C71").Value = 0.33
$ws.Range("E71").Value = 10
$ws.Range("G71").Value = "1x a cada 6 meses - irregular (preferencialmente na 2ª quinzena)"
$ws.Range("H71").Value = (Get-Date -Year 2025 -Month 6 -Day 2 -Hour 19 -Minute 9 -Second 6)
$ws.Range("I71").Value = (Get-Date -Year 2025 -Month 12 -Day 2 -Hour 19 -Minute 9 -Second 6)
$ws.Range("I71").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("J71").Value = "ATIVO"

# --- Row 111 (id_cliente 28458, BEMOL S/A): refreshed purchase stats ---
$ws.Range("B111").Value = 0.92
$ws.Range("C111").Value = 0.83
$ws.Range("E111").Value = 15123
$ws.Range("H111").Value = (Get-Date -Year 2025 -Month 6 -Day 2 -Hour 17 -Minute 25 -Second 22)
$ws.Range("I111").Value = (Get-Date -Year 2025 -Month 6 -Day 3 -Hour 17 -Minute 25 -Second 22)

"edit applied"
